# add alert to the smart contract
# Patient PA001 sheet: shift the observation window forward by one day
# (Observation 2 now carries Observation 1's timestamp/values) and clear
# out what used to be Observation 3 (column F) so the alert can be
# recomputed against the fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Date
$ws.Cells.Item(3, 4).Value = 44334.963530092595   # D3
$ws.Cells.Item(3, 5).Value = 44334.963530092595   # E3
$ws.Cells.Item(3, 6).ClearContents()              # F3

# Row 4 - W (weight)
$ws.Cells.Item(4, 5).Value = 30                   # E4
$ws.Cells.Item(4, 6).ClearContents()              # F4
$ws.Cells.Item(4, 6).NumberFormat = "General"

# Row 5 - BS
$ws.Cells.Item(5, 5).Value = 1.4                  # E5
$ws.Cells.Item(5, 6).ClearContents()              # F5
$ws.Cells.Item(5, 6).NumberFormat = "General"

# Row 6 - SE
$ws.Cells.Item(6, 5).Value = 0.2                  # E6
$ws.Cells.Item(6, 6).ClearContents()              # F6

# Row 7 - SSE
$ws.Cells.Item(7, 5).Value = 1                    # E7
$ws.Cells.Item(7, 6).ClearContents()              # F7

# Row 8 - Age
$ws.Cells.Item(8, 6).ClearContents()              # F8
$ws.Cells.Item(8, 6).NumberFormat = "General"

# Row 9 - counter
$ws.Cells.Item(9, 2).Value = 2                    # B9
